$d = $word.ActiveDocument

# --- First paragraph: the hidden "**ID__...__ID**" bookmark-id paragraph ---
$p1 = $d.Paragraphs(1)

# Add a (space-only, no visible line) paragraph border around paragraph 1
$borders = $p1.Range.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Drop the trailing space run that followed the id text, then update
# the id itself to the new topic identifier ("last minute updates").
$d.Content.Find.Execute("**ID__AFFARS_pgi_5315_topic_10__ID** ", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "**ID__AFFARS_pgi_5315_topic_10__ID**", 2) | Out-Null

$d.Content.Find.Execute("**ID__AFFARS_pgi_5315_topic_10__ID**", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "**ID__AFFARS_AF_PGI_5315_408_90__ID**", 2) | Out-Null
